# Add a new worksheet "Sheet2" right after "Sheet1", populate it with the
# two new city names, autofit column A, select B3 on it, and leave it as
# the active sheet/tab (matches workbookView activeTab="1" plus Sheet2's
# tabSelected/selection in the target workbook).
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Kolkata"
$ws2.Range("A2").Value = "Hyderabad"

$ws2.Columns("A:A").AutoFit() | Out-Null
$ws2.Range("B3").Select() | Out-Null
